$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right after
#    the title (Heading1) paragraph.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Build a new bold paragraph containing the page title text
#    ("Play Candy Witch for Free - Review of Gameplay, Bonuses & Payouts").
#    We create it right after a plain paragraph (no paragraph style, no
#    inherited run formatting) so the new paragraph itself starts out clean,
#    then we move it into position just before the final paragraph.
# ---------------------------------------------------------------------------
$plainSourcePara = $d.Paragraphs.Item(3)
$plainSourcePara.Range.InsertParagraphAfter()

$titleText = "Play Candy Witch for Free - Review of Gameplay, Bonuses & Payouts"
$builtPara = $d.Paragraphs.Item(4)
$builtPara.Range.Text = $titleText

$builtPara = $d.Paragraphs.Item(4)
$boldRange = $d.Range($builtPara.Range.Start, $builtPara.Range.Start + $titleText.Length)
$boldRange.Font.Bold = $true

# Move the freshly-built paragraph (cut) to just before the last paragraph in
# the document (the one that currently holds the image-prompt text).
$builtPara = $d.Paragraphs.Item(4)
$builtPara.Range.Cut()

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$destination = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$destination.Paste()

# ---------------------------------------------------------------------------
# 3) Replace the text of the (now) final paragraph - the italic image prompt -
#    with the meta-description copy, keeping its italic formatting.
# ---------------------------------------------------------------------------
$oldText = "Create a cartoon-style feature image for Candy Witch that features a happy Maya warrior with glasses. The image should have a bright and colorful background, with the Maya warrior positioned in the center of the frame. He should be smiling and holding a handful of candy in one hand, with the other hand raised up in a magic spell-casting pose. His glasses should be oversized and cartoonish, with a reflection of the glow from the enchanted forest in the lenses. In the background, there should be hints of the forest and the moonlight. The overall image should convey a sense of fun and whimsy while also capturing the magic and excitement of the Candy Witch slot game."
$newText = "Try Candy Witch for free and discover two exciting bonus games with Sticky Wilds and progressive multipliers, plus great payouts."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
